# ============================================================
# Edit script: Added Arm Servo info comments and torque calculations
# ============================================================
$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: remove the old "_GoBack" bookmark that sat right after
# "Replaced " in the eye-LED wire bullet point.
# ------------------------------------------------------------------
try {
    $oldBm = $d.Bookmarks.Item("_GoBack")
    $oldBm.Delete()
} catch {
}

# ------------------------------------------------------------------
# Change 2: append the new "arm servo torque" commentary right after
# the "...extension cords ... disconnect inside." sentence, then move
# the "_GoBack" bookmark to the end of the freshly added text.
# ------------------------------------------------------------------
$rFind = $d.Content
$found = $rFind.Find.Execute("extension cords for each servo or group of servos that can disconnect inside.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertAt = $rFind.End
$insertRange = $d.Range($insertAt, $insertAt)
$bigXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:br/></w:r><w:r><w:t xml:space="preserve">So I''m attempting it without using any torque multiplying servo </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>gearbox''s</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> at first because I want the movement as fast as possible.   For the shoulder I plan on using two relatively in expensive </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>GoBILDA</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">  2000</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> Series Dual Mode Analog Servo (25-2) (300 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>oz</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>in  stall</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> torque at 6Volts with no load servo speed of 0.2 sec/60degrees, weight 2.12 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>oz</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">)  for roll and yaw and a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Savox</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> SA1230SG Coreless Digital Servo (499.9 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>oz</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">-in at 6 Volts,  with no load speed 0.16 sec/60 degree,  weight 2.8 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>oz</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">) for pitch (i.e. lifting the arm up and down).   Then another </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>GoBILDA</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">  2000</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> Series Dual Mode Analog Servo (25-2) for the elbow.  My skeleton arm weighs 6.4 </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>oz</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">  and</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> is 24" long with its balance point around the 12" mark from the shoulder.  So based on that I believe the bare min torque to hold the arm still straight out perpendicular to the body is (6.4 + 2.12) </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>oz</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> * 12"  =  102 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>oz</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">-in of torque (assuming the elbow motor is right at the elbow although I could put it at the shoulder and use a light weight linkage to the elbow).    So I figure a stall torque of say 3 times that or more should be good enough to move the arm at a reasonable rate and decelerate it ok and not overstress the servo.   I suppose I could calculate the speed assuming the motor stall </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">torque and arm weight distribution but I''m guessing 3 times the bare min still arm torque is good enough .  My 499 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>oz</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-in motor is basically 5 times the bare min still arm torque.  Most of the time the arm will be in the resting positing hanging straight down.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertRange.InsertXML($bigXml)
# InsertXML minted a fresh paragraph mark right at $insertAt; delete it
# so the new runs stay inside the original list-item paragraph.
$mergeMark = $d.Range($insertAt, $insertAt + 1)
$mergeMark.Delete()

# Re-anchor the "_GoBack" bookmark at the very end of the new text.
$rTail = $d.Content
$foundTail = $rTail.Find.Execute("Most of the time the arm will be in the resting positing hanging straight down.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tailEnd = $rTail.End
$bmRange = $d.Range($tailEnd, $tailEnd + 1)
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# Change 3: split "...as it would be easier to control " into three
# runs so "would" gets wrapped in proofErr gramStart/gramEnd markers.
# ------------------------------------------------------------------
$r3 = $d.Content
$found3 = $r3.Find.Execute("air cylinder would work for back up/down as it would be easier to control ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r3b = $d.Range($r3.Start, $r3.End)
$xml3 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">air cylinder would work for back up/down as it </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>would</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> be easier to control </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r3b.InsertXML($xml3)

# ------------------------------------------------------------------
# Change 4: merge the two runs around the old lastRenderedPageBreak
# back into a single run (removing the stray page-break marker).
# ------------------------------------------------------------------
$r4 = $d.Content
$found4 = $r4.Find.Execute("in some manner instead of connecting off LEDs on that audio level board for audio level.  The audio needs to have enough release to hold the level so that the Arduino is not taxed sampling that audio level.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r4b = $d.Range($r4.Start, $r4.End)
$xml4 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>in some manner instead of connecting off LEDs on that audio level board for audio level.  The audio needs to have enough release to hold the level so that the Arduino is not taxed sampling that audio level.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r4b.InsertXML($xml4)

Write-Host "Done. found1-4: $found $found3 $found4"
